# Atualização automática de preços de eletricidade
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45929

$ws.Range("B2").Value = 51.6
$ws.Range("C2").Value = 47.7
$ws.Range("D2").Value = 45
$ws.Range("E2").Value = 32
$ws.Range("F2").Value = 30
$ws.Range("G2").Value = 32.5
$ws.Range("H2").Value = 50
$ws.Range("I2").Value = 74.40000000000001
$ws.Range("J2").Value = 77.8
$ws.Range("K2").Value = 73.08
$ws.Range("L2").Value = 59.33
$ws.Range("M2").Value = 40.14
$ws.Range("N2").Value = 35
$ws.Range("O2").Value = 28.3
$ws.Range("P2").Value = 21.99
$ws.Range("Q2").Value = 27.94
$ws.Range("R2").Value = 35
$ws.Range("S2").Value = 52
$ws.Range("T2").Value = 85
$ws.Range("U2").Value = 101.13
$ws.Range("V2").Value = 133.26
$ws.Range("W2").Value = 105.68
$ws.Range("X2").Value = 90.04000000000001
$ws.Range("Y2").Value = 88.59999999999999
$ws.Range("Z2").Value = 59.06

$ws.Range("AB2").Value = 104.4
$ws.Range("AD2").Value = 119.47
$ws.Range("AF2").Value = 93.06

$ws.Range("AG2").Value = "0h-17h"
